$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "maa://24702 (94.24), maa://25390 (95.92), maa://36681 (87.01)"
$ws.Range("L2").Value = "*maa://24633 (55.48), *maa://30515 (69.0), *maa://34787 (72.06), ***maa://20792 (11.93), maa://39402 (87.5), ***maa://29083 (27.78)"
$ws.Range("L3").Value = "*maa://22880 (66.12), maa://20276 (84.81), *maa://22749 (72.73)"
$ws.Range("T3").Value = "maa://24617 (88.89), **maa://20790 (44.78), ***maa://37170 (20.0)"
$ws.Range("X3").Value = "maa://27396 (84.42), maa://27484 (96.12), maa://27480 (82.86)"
$ws.Range("D4").Value = "maa://24632 (93.24), **maa://24303 (33.33), maa://22499 (86.67), maa://22746 (100.0)"
$ws.Range("D7").Value = "maa://21955 (94.12)"
$ws.Range("A8").Value = "更新日期：2024.12.22 13:16:54"
$ws.Range("D8").Value = "*maa://21476 (72.34), **maa://39431 (45.45), *maa://37551 (57.14)"
$ws.Range("P8").Value = "maa://32931 (84.62), *maa://21916 (61.29), maa://23252 (92.42), maa://37496 (96.3), **maa://22759 (45.45)"
$ws.Range("X8").Value = "maa://21411 (95.91)"
$ws.Range("P9").Value = "maa://22736 (81.32)"
$ws.Range("X9").Value = "maa://26223 (97.46)"
$ws.Range("AF9").Value = "maa://26206 (89.32), *maa://22865 (50.98)"
$ws.Range("D10").Value = "***maa://25695 (19.44), **maa://32237 (40.48), ***maa://34206 (21.74), ***maa://39951 (17.07), ***maa://39243 (28.57)"
$ws.Range("AB12").Value = "maa://23669 (95.39), maa://36677 (93.88), maa://39872 (90.0)"
$ws.Range("AF12").Value = "*maa://28932 (78.79), *maa://20106 (63.96), *maa://22769 (64.29)"
$ws.Range("X13").Value = "*maa://34957 (79.37), *maa://22768 (51.61)"
$ws.Range("L14").Value = "maa://26245 (96.5), maa://21288 (96.3), maa://39841 (95.06), maa://36682 (97.37)"
$ws.Range("D15").Value = "*maa://22743 (76.96), maa://22734 (84.03), *maa://30808 (64.52), **maa://36048 (33.33)"
$ws.Range("P15").Value = "maa://24762 (90.26), *maa://22727 (70.0)"
$ws.Range("AF15").Value = "maa://21364 (80.91), *maa://22766 (70.91), *maa://36666 (78.57)"
$ws.Range("H18").Value = "maa://24421 (90.34)"
$ws.Range("T19").Value = "maa://24386 (99.04)"
$ws.Range("AB19").Value = "*maa://30709 (63.14), *maa://36668 (55.84)"
$ws.Range("H20").Value = "maa://22864 (89.26)"
$ws.Range("AF21").Value = "maa://22524 (94.63), *maa://22432 (76.67)"
$ws.Range("L22").Value = "maa://27127 (85.29), *maa://22751 (73.85)"
$ws.Range("X22").Value = "maa://21282 (98.46), *maa://37649 (69.23)"
$ws.Range("L23").Value = "maa://39756 (94.05), maa://39875 (93.75)"
$ws.Range("P23").Value = "maa://30587 (91.98), *maa://29748 (75.78), ***maa://29785 (16.42), *maa://37566 (71.43)"
$ws.Range("X24").Value = "maa://29988 (86.73), maa://23504 (93.28), **maa://22892 (39.58), *maa://25141 (77.6), maa://36663 (80.6), ***maa://22815 (23.08)"
$ws.Range("H25").Value = "*maa://29063 (73.65), *maa://25311 (74.0), ***maa://22725 (4.84)"
$ws.Range("AF25").Value = "maa://20108 (96.27), maa://24621 (96.58), maa://36676 (96.67), maa://22771 (85.71), maa://37772 (100.0)"
$ws.Range("AB26").Value = "maa://42235 (93.24)"
$ws.Range("T28").Value = "maa://23263 (94.95), *maa://29765 (60.53)"
$ws.Range("X28").Value = "maa://39929 (89.74), ***maa://39723 (14.29), maa://41749 (90.74)"
$ws.Range("L29").Value = "maa://28432 (92.86), *maa://28440 (76.6), maa://31400 (100.0), *maa://28650 (71.43)"
$ws.Range("AF29").Value = "*maa://24080 (69.05), ***maa://34960 (8.33), *maa://42865 (78.79)"
$ws.Range("L31").Value = "maa://35926 (93.98), maa://36258 (83.16), *maa://43904 (71.43)"
$ws.Range("H39").Value = "maa://25199 (84.82), maa://36670 (87.21), maa://30434 (89.23), ***maa://25036 (16.0), **maa://44165 (50.0)"
$ws.Range("P41").Value = "**maa://35616 (38.24), maa://43177 (85.71)"
$ws.Range("H45").Value = "maa://21229 (84.86), maa://30807 (95.45), *maa://22767 (55.0), ***maa://20796 (13.79), *maa://42459 (75.0)"
$ws.Range("H47").Value = "maa://27410 (96.18), maa://29661 (97.84), maa://28038 (84.62)"
